# Update the header for column F (H1: period changed from "Январь - Июнь" to "Январь - Сентябрь")
# and refresh all the per-region data in column F with the new period's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "2023`n Январь - Сентябрь"

$ws.Range("F2").Value  = 3802.6
$ws.Range("F3").Value  = 3169.4
$ws.Range("F4").Value  = 2921.4
$ws.Range("F5").Value  = 3385.6
$ws.Range("F6").Value  = 2041.4
$ws.Range("F7").Value  = 13887.1
$ws.Range("F8").Value  = 4384.9
$ws.Range("F9").Value  = 1563.2
$ws.Range("F10").Value = 1728.9
$ws.Range("F11").Value = 4509.7
$ws.Range("F12").Value = 3261.7
$ws.Range("F13").Value = 2238.4
$ws.Range("F14").Value = 4548.8
$ws.Range("F15").Value = 3374.5
$ws.Range("F16").Value = 2884.7
$ws.Range("F17").Value = 1268.1
$ws.Range("F18").Value = 5537.3
$ws.Range("F19").Value = 4283.6
$ws.Range("F20").Value = 5673.3
$ws.Range("F21").Value = 6673.5
$ws.Range("F22").Value = 2033.5
